# Update 南宁-漫展信息.xlsx with refreshed "想去人数" (and a couple "最低票价")
# numbers for the gh-pages data snapshot generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 5380
$ws.Range("F4").Value  = 11397
$ws.Range("G4").Value  = 62
$ws.Range("F5").Value  = 279
$ws.Range("F6").Value  = 587
$ws.Range("F7").Value  = 167
$ws.Range("F8").Value  = 249
$ws.Range("F9").Value  = 978
$ws.Range("F10").Value = 95

# --- Sheet "演出" -------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 10

# --- Sheet "全部类型" ---------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 5380
$ws.Range("F7").Value  = 11397
$ws.Range("G7").Value  = 62
$ws.Range("F8").Value  = 279
$ws.Range("F9").Value  = 587
$ws.Range("F10").Value = 167
$ws.Range("F11").Value = 10
$ws.Range("F13").Value = 249
$ws.Range("F14").Value = 978
$ws.Range("F16").Value = 95
